$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.20%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.72%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.562'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.41%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08083'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.50%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.918'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.45%'
$ws.Range("B7").Value = 'BTSEToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '2.574'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-7.31%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9512'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.70%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1187'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.57%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1850'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.07%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09777'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.13%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04494'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '4.39%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1067'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.17%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001283'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.29%'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04191'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-4.53%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005850'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.69%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.387'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-4.77%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.295'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.78%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3463'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.55%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.25'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '16.38%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1419'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.31%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001244'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.17%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004368'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.93%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001189'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.09%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.98%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02686'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '1.37%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05555'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.41%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007557'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.63%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1408'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.49%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.008312'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-14.89%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002014'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.68%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008893'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-7.52%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007165'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.36%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.84%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.003154'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-8.99%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002269'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-0.78%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.84%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.84%'
